$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$shp = $ws.Shapes.Item(1)
$shp.Left = 100
$shp.Top = 120
$shp.Width = 150
$shp.Height = 75
Write-Host "Set Top/Left/Width/Height OK"
